$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shared string / header for column D ("ibr")
$ws.Range("D1").Value = "ibr"

# Fill D699:D800 with the new "ibr" series values
$dData = @(0.049868095238095239, 0.046889500000000001, 0.047448181818181813, 0.046769499999999999, 0.044384736842105268, 0.040678571428571432, 0.038317499999999997, 0.035555555555555549, 0.031535909090909094, 0.032013333333333331, 0.032252777777777776, 0.032142272727272728, 0.032627999999999997, 0.032064285714285712, 0.032530909090909089, 0.032496842105263164, 0.032312500000000001, 0.032483809523809519, 0.032688000000000002, 0.032195999999999995, 0.0326515, 0.036377142857142858, 0.039267777777777783, 0.040868260869565215, 0.042753684210526312, 0.044670000000000015, 0.044844999999999996, 0.044998333333333335, 0.044978999999999998, 0.045058000000000001, 0.045023000000000001, 0.04498428571428572, 0.045145499999999998, 0.045115789473684213, 0.045102631578947373, 0.045195909090909092, 0.045866842105263164, 0.04632, 0.049519523809523812, 0.054362631578947364, 0.057374000000000001, 0.059357894736842109, 0.062585238095238097, 0.064151, 0.066457142857142854, 0.071649500000000005, 0.074389999999999998, 0.076189999999999994, 0.078075454545454545, 0.077508181818181809, 0.077477000000000004, 0.077370999999999995, 0.076184500000000002, 0.07399571428571429, 0.074383999999999992, 0.070619545454545463, 0.067951666666666674, 0.062965238095238088, 0.060464000000000004, 0.055301052631578945, 0.053307619047619041, 0.052057619047619047, 0.051895714285714281, 0.048726000000000005, 0.047236111111111118, 0.046452380952380953, 0.044981999999999994, 0.044464736842105265, 0.043272380952380951, 0.042439047619047619, 0.042483684210526312, 0.042512499999999995, 0.042515238095238106, 0.042504999999999994, 0.042505909090909094, 0.042521500000000004, 0.042494736842105273, 0.042499047619047617, 0.042493499999999997, 0.042485000000000002, 0.042495999999999999, 0.042496363636363638, 0.042481666666666668, 0.042450454545454548, 0.042546999999999995, 0.042555714285714287, 0.042531363636363631, 0.042523684210526318, 0.042507500000000004, 0.04252809523809524, 0.042522000000000004, 0.04136571428571429, 0.033949, 0.02967052631578947, 0.02568368421052631, 0.02342909090909091, 0.02117157894736842, 0.019235909090909088, 0.017491428571428571, 0.01752894736842105, 0.017519, 0.017234736842105261)

for ($i = 0; $i -lt $dData.Length; $i++) {
    $ws.Cells.Item(699 + $i, 4).Value = $dData[$i]
}

# Update the active selection to match the committed state
$ws.Range("E699").Select()

